$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 208, pushing existing rows 208..300 down to 209..301
$ws.Rows.Item(208).Insert()

# Populate the newly inserted row 208 with its data
$ws.Range("A208").Value = 5
$ws.Range("B208").Value = "Macroferia Regional de Talca"
$ws.Range("C208").Value = "Maule"
$ws.Range("D208").Value = 44636
$ws.Range("E208").Value = 7
$ws.Range("F208").Value = 100114013
$ws.Range("G208").Value = "Zanahoria"
$ws.Range("H208").Value = "Sin especificar"
$ws.Range("I208").Value = "Primera"
$ws.Range("J208").Value = 300
$ws.Range("K208").Value = 7000
$ws.Range("L208").Value = 7000
$ws.Range("M208").Value = 7000
$ws.Range("N208").Value = "$/saco 20 kilos"
$ws.Range("O208").Value = "Región de Ñuble"
$ws.Range("P208").Value = 350
$ws.Range("Q208").Value = 20
$ws.Range("R208").Value = "Hortaliza"
